$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quantities for the two new rows (plain numbers, not shared strings)
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1

# Link for the existing "Stirnlampe" item (row 6 / Item_Nr 5)
$ws.Range("D6").Value = "https://www.bauhaus.info/stirnlampen/profi-depot-stirnlampe/p/25568884"

# New item (row 7 / Item_Nr 6): link filled in first
$ws.Range("D7").Value = "https://www.amazon.de/s?k=oneplus+6t+panzerglas&__mk_de_DE=%C3%85M%C3%85%C5%BD%C3%95%C3%91&crid=2O2D66FWAGP4K&sprefix=oneplus+6t+%2Caps%2C190&ref=nb_sb_ss_i_3_11"

# New item (row 8 / Item_Nr 7): name filled in before row 7's name
$ws.Range("B8").Value = "Rhinoshield solidsuit schutzh$([char]0x00FC)lle oneplus 6t"

# New item (row 7 / Item_Nr 6): name
$ws.Range("B7").Value = "Schutzfolie/Panzerglas Oneplus 6t"

# New item (row 8 / Item_Nr 7): link filled in last
$ws.Range("D8").Value = "https://rhinoshield.de/pages/shop/OnePlus?device=oneplus-6t&category=solidsuit&collection=solidsuit-android"

# Column B needs to auto-fit its width to the new, longer text
$ws.Columns.Item(2).ColumnWidth = 34.3

# Move active selection to D8, matching the editor's last-touched cell
$ws.Range("D8").Select()
